$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To fix")

# Insert a new header row at the top
$ws.Rows.Item(1).Insert()

# Header row content
$ws.Range("A1").Value = "Bug/Feature"
$ws.Range("B1").Value = "Status"

# Status values for each bug row (rows 2-11)
$ws.Range("B2").Value = "Done"
$ws.Range("B3").Value = "Done"
$ws.Range("B4").Value = "Open"
$ws.Range("B5").Value = "Open"
$ws.Range("B6").Value = "Rejected"
$ws.Range("B7").Value = "Open"
$ws.Range("B8").Value = "Done"
$ws.Range("B9").Value = "?"
$ws.Range("B10").Value = "?"
$ws.Range("B11").Value = "Open"

# Column width (bestFit-like)
$ws.Columns.Item(1).ColumnWidth = 84.45182291666667

# Selection matches target
$ws.Range("A14").Select()

# Mark header row with explicit (no-op) fill so a distinct style is recorded
$ws.Range("A1:B1").Interior.ColorIndex = -4142

# AutoFilter on header row (also creates the hidden _FilterDatabase defined name)
$ws.Range("A1:B1").AutoFilter()
$n = $ws.Names.Add("_xlnm._FilterDatabase", "='To fix'!`$A`$1:`$B`$1")
$n.Visible = $false
